# Insert 3 new weekly price records before the current row 312, pushing the
# existing rows 312:397 down to 315:400 (dimension grows from R397 to R400).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("312:314").Insert()

# Row 312: Repollo - Copenhague - Primera
$ws.Cells.Item(312, 1).Value = 4
$ws.Cells.Item(312, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(312, 3).Value = "Los Lagos"
$ws.Cells.Item(312, 4).Value = 44642
$ws.Cells.Item(312, 5).Value = 10
$ws.Cells.Item(312, 6).Value = 100112006
$ws.Cells.Item(312, 7).Value = "Repollo"
$ws.Cells.Item(312, 8).Value = "Copenhague"
$ws.Cells.Item(312, 9).Value = "Primera"
$ws.Cells.Item(312, 10).Value = 400
$ws.Cells.Item(312, 11).Value = 2000
$ws.Cells.Item(312, 12).Value = 2000
$ws.Cells.Item(312, 13).Value = 2000
$ws.Cells.Item(312, 14).Value = "`$/unidad"
$ws.Cells.Item(312, 15).Value = "Región Metropolitana"
$ws.Cells.Item(312, 16).Value = 2000
$ws.Cells.Item(312, 17).Value = 1
$ws.Cells.Item(312, 18).Value = "Hortaliza"

# Row 313: Repollo - Copenhague - Segunda
$ws.Cells.Item(313, 1).Value = 4
$ws.Cells.Item(313, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(313, 3).Value = "Los Lagos"
$ws.Cells.Item(313, 4).Value = 44642
$ws.Cells.Item(313, 5).Value = 10
$ws.Cells.Item(313, 6).Value = 100112006
$ws.Cells.Item(313, 7).Value = "Repollo"
$ws.Cells.Item(313, 8).Value = "Copenhague"
$ws.Cells.Item(313, 9).Value = "Segunda"
$ws.Cells.Item(313, 10).Value = 400
$ws.Cells.Item(313, 11).Value = 1800
$ws.Cells.Item(313, 12).Value = 1800
$ws.Cells.Item(313, 13).Value = 1800
$ws.Cells.Item(313, 14).Value = "`$/unidad"
$ws.Cells.Item(313, 15).Value = "Región Metropolitana"
$ws.Cells.Item(313, 16).Value = 1800
$ws.Cells.Item(313, 17).Value = 1
$ws.Cells.Item(313, 18).Value = "Hortaliza"

# Row 314: Repollo - Crespo record - Primera
$ws.Cells.Item(314, 1).Value = 4
$ws.Cells.Item(314, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(314, 3).Value = "Los Lagos"
$ws.Cells.Item(314, 4).Value = 44642
$ws.Cells.Item(314, 5).Value = 10
$ws.Cells.Item(314, 6).Value = 100112006
$ws.Cells.Item(314, 7).Value = "Repollo"
$ws.Cells.Item(314, 8).Value = "Crespo record"
$ws.Cells.Item(314, 9).Value = "Primera"
$ws.Cells.Item(314, 10).Value = 800
$ws.Cells.Item(314, 11).Value = 1700
$ws.Cells.Item(314, 12).Value = 1700
$ws.Cells.Item(314, 13).Value = 1700
$ws.Cells.Item(314, 14).Value = "`$/unidad"
$ws.Cells.Item(314, 15).Value = "Región Metropolitana"
$ws.Cells.Item(314, 16).Value = 1700
$ws.Cells.Item(314, 17).Value = 1
$ws.Cells.Item(314, 18).Value = "Hortaliza"
